# Remove the scraped image-path and HTML description text that used to
# live in columns E ("Image") and F ("Description") for a large batch of
# product rows. The cell stays in place (still t="inlineStr") but is left
# with no value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UpdatedList")

# Rows whose Image (E) AND Description (F) text must both be cleared.
$ws.Range("E3:F4").ClearContents()
$ws.Range("E6:F14").ClearContents()
$ws.Range("E18:F18").ClearContents()
$ws.Range("E24:F28").ClearContents()
$ws.Range("E30:F35").ClearContents()
$ws.Range("E53:F54").ClearContents()

# Rows whose Image (E) text must be cleared (Description was already empty).
$ws.Range("E20:E23").ClearContents()
$ws.Range("E37:E41").ClearContents()
$ws.Range("E43:E45").ClearContents()
$ws.Range("E55:E62").ClearContents()
$ws.Range("E64").ClearContents()
$ws.Range("E66:E72").ClearContents()
